# Regenerate the "K" column (column G) values for save_data sheet.
# These values were recomputed (std/mean, calc and write s_vals) and
# replace the previous "Strike#" based figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 4
    4  = 3
    5  = 0
    6  = 5
    7  = 1
    8  = 4
    9  = 3
    10 = 2
    11 = 3
    12 = 3
    13 = 2
    14 = 2
    15 = 4
    16 = 1
    17 = 5
    18 = 4
    19 = 2
    20 = 2
    21 = 6
    22 = 3
    23 = 1
    24 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
